# "Add cantrals by cantons"
#
# The sheet used to have a two-row header (row1 + row2) above 8 rows of
# hydro-plant data. This rewrites it to a single-row header with more
# descriptive column titles, adds two new leading id columns (idx/idx2
# -> existing A/B data columns), and shifts the data block up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old row 2 ("Hiver"/"Ete"/"Annee" sub-header) goes away entirely;
# deleting it shifts all the data rows (old 3..10) up to (2..9).
$ws.Rows.Item(2).Delete()

# Columns A:E of the header keep the default (unstyled) look.
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

# Columns F:K keep the same (Arial 9pt) look the data cells use, with new
# more descriptive header text.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Match the workbook's new saved selection (first data row).
$ws.Range("A2:K2").Select() | Out-Null
